$d = $word.ActiveDocument

# 1. Remove the "Phone" row from the Personal info table (first table in the
#    document). The row contains "Phone" / "+48 501 196 255". (Cell.Range.Text
#    carries trailing cell-end marker characters, so match with a wildcard.)
$table = $d.Tables.Item(1)
for ($i = $table.Rows.Count; $i -ge 1; $i--) {
    $row = $table.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text -like "Phone*") {
        $row.Delete()
    }
}

# 2. Update the Clear2Pay Poland employment dates (it is no longer a current
#    position, it now has an end date).
$d.Content.Find.Execute("since 05.2016: Clear2Pay Poland", $true, $false,
                         $false, $false, $false, $true, 1, $false,
                         "05.2016 - 07.2016: Clear2Pay Poland", 2)

# 3. Tweak the wording of one of the Clear2Pay bullet points.
$d.Content.Find.Execute("Development of internal Java based tools", $true,
                         $false, $false, $false, $false, $true, 1, $false,
                         "Development of internal, Java based tool", 2)
